# Append the latest gold-price row (10-12-2025) to the data sheet,
# mirroring the existing rows' layout: col A = date text, col B = price text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$dateText = "10-12-2025"
$priceText = "The price of gold in India today is ₹13,031 per gram for 24 karat gold, ₹11,945 per gram for 22 karat gold and ₹9,773 per gram for 18 karat gold (also called 999 gold)."

# Force text storage so the date-like string isn't auto-converted to a
# date serial number, then write the values.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = $dateText
$ws.Cells.Item($newRow, 2).Value = $priceText

# Re-apply the same cell formatting (border/fill/wrap) used by the rest
# of the table, copying it from the row directly above.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)

$ws.Cells.Item($newRow - 1, 2).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4122)

$excel.CutCopyMode = $false
